$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.965.57"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.848.49"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'309.66"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.4765"
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'0.07219"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "'0.9267"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").Value = "'19.68"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'0.07712"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.803.80"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'5.319"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "'6.413"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'88.83"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'0.000008634"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "27.001.44"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'10.65"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "'1.935"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'152.58"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'1.996"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'114.18"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'4.951"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "'0.08875"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'3.320"
$ws.Range("E31").Value = "  +5.39%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'0.7431"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "'4.493"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'2.714"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'1.113"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'0.05261"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("D39").Value = "'2.982"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "'0.5188"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "'6.985"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "'8.191"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("D44").Value = "'10.62"
$ws.Range("E44").Value = "  +6.19%  "
$ws.Range("D45").Value = "'0.4725"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'101.57"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("D48").Value = "'1.602"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'66.10"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").Value = "'0.06028"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.8879"
$ws.Range("E51").Value = "  +4.28%  "
